$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing values (rounded figures)
$ws.Range("D22").Value = 0.224211961
$ws.Range("C23").Value = -0.18190582

# Add new D23 value
$ws.Range("D23").Value = 0.6740699019999999

# Add new row 24
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A24").Value = "2025-09-04_diff"
$ws.Range("B24").Value = -0.431278794
$ws.Range("C24").Value = 0.118688665
